{"js": "// Apply text replacements described by the diff:\n// date header + 25 two-digit multiplication expressions in the table.\nconst replacements = [\n  [\"2025-03-12 Wednesday\", \"2025-03-13 Thursday\"],\n  [\"22\u00d729=\", \"71\u00d725=\"],\n  [\"51\u00d754=\", \"35\u00d723=\"],\n  [\"22\u00d771=\", \"72\u00d736=\"],\n  [\"15\u00d775=\", \"43\u00d722=\"],\n  [\"49\u00d775=\", \"20\u00d782=\"],\n  [\"46\u00d779=\", \"78\u00d714=\"],\n  [\"34\u00d791=\", \"54\u00d712=\"],\n  [\"57\u00d777=\", \"84\u00d797=\"],\n  [\"58\u00d759=\", \"46\u00d769=\"],\n  [\"87\u00d772=\", \"96\u00d731=\"],\n  [\"36\u00d776=\", \"24\u00d720=\"],\n  [\"37\u00d716=\", \"60\u00d722=\"],\n  [\"80\u00d779=\", \"22\u00d722=\"],\n  [\"51\u00d714=\", \"62\u00d729=\"],\n  [\"82\u00d793=\", \"83\u00d757=\"],\n  [\"86\u00d737=\", \"66\u00d787=\"],\n  [\"20\u00d737=\", \"43\u00d785=\"],\n  [\"87\u00d763=\", \"84\u00d724=\"],\n  [\"43\u00d774=\", \"57\u00d779=\"],\n  [\"62\u00d712=\", \"31\u00d751=\"],\n  [\"36\u00d753=\", \"18\u00d789=\"],\n  [\"23\u00d795=\", \"49\u00d742=\"],\n  [\"87\u00d783=\", \"74\u00d714=\"],\n  [\"34\u00d794=\", \"50\u00d770=\"],\n  [\"88\u00d718=\", \"30\u00d774=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply text replacements described by the diff:\n# date header + 25 two-digit multiplication expressions in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-03-12 Wednesday\", \"2025-03-13 Thursday\"),\n  @(\"22\u00d729=\", \"71\u00d725=\"),\n  @(\"51\u00d754=\", \"35\u00d723=\"),\n  @(\"22\u00d771=\", \"72\u00d736=\"),\n  @(\"15\u00d775=\", \"43\u00d722=\"),\n  @(\"49\u00d775=\", \"20\u00d782=\"),\n  @(\"46\u00d779=\", \"78\u00d714=\"),\n  @(\"34\u00d791=\", \"54\u00d712=\"),\n  @(\"57\u00d777=\", \"84\u00d797=\"),\n  @(\"58\u00d759=\", \"46\u00d769=\"),\n  @(\"87\u00d772=\", \"96\u00d731=\"),\n  @(\"36\u00d776=\", \"24\u00d720=\"),\n  @(\"37\u00d716=\", \"60\u00d722=\"),\n  @(\"80\u00d779=\", \"22\u00d722=\"),\n  @(\"51\u00d714=\", \"62\u00d729=\"),\n  @(\"82\u00d793=\", \"83\u00d757=\"),\n  @(\"86\u00d737=\", \"66\u00d787=\"),\n  @(\"20\u00d737=\", \"43\u00d785=\"),\n  @(\"87\u00d763=\", \"84\u00d724=\"),\n  @(\"43\u00d774=\", \"57\u00d779=\"),\n  @(\"62\u00d712=\", \"31\u00d751=\"),\n  @(\"36\u00d753=\", \"18\u00d789=\"),\n  @(\"23\u00d795=\", \"49\u00d742=\"),\n  @(\"87\u00d783=\", \"74\u00d714=\"),\n  @(\"34\u00d794=\", \"50\u00d770=\"),\n  @(\"88\u00d718=\", \"30\u00d774=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $result = $find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n  )\n\n  if (-not $result) {\n    throw \"Replacement failed for: $oldText\"\n  }\n}\n"}
